# Feature: new function: Write to baseworkspace
#
# Adds a new row (IndexNum 37) to the "表1" Excel Table on Sheet1 for the
# newly introduced "Write result to BaseWorkSpace" text entry, in English,
# Chinese and Japanese - growing the table from A1:E38 to A1:E39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row; ListRows.Add() appends a row at the bottom and
# expands the table/autofilter ref (A1:E38 -> A1:E39) and sheet dimension.
$newRow = $lo.ListRows.Add()

# Worksheet row the freshly-added table row landed on (39).
$dataRow = $newRow.Range.Row

$ws.Cells.Item($dataRow, 1).Value = 37
$ws.Cells.Item($dataRow, 2).Value = "Write result to BaseWorkSpace"
$ws.Cells.Item($dataRow, 3).Value = "合并结果写到基础工作区"
$ws.Cells.Item($dataRow, 4).Value = "Write result to BaseWorkSpace"
$ws.Cells.Item($dataRow, 5).Value = "マージ結果をベースワークスペースに書き込む"

# Match the wrap/vertical-center formatting used by the other
# Chinese/English/Japanese text columns, and the taller row height that
# goes with the wrapped, two-line Chinese/Japanese translations.
$textRange = $ws.Range($ws.Cells.Item($dataRow, 3), $ws.Cells.Item($dataRow, 5))
$textRange.WrapText = $true
$textRange.VerticalAlignment = -4108
$ws.Rows.Item($dataRow).RowHeight = 28

# Reflect the new scroll position/selection left behind by the edit.
$ws.Range("E42").Select()
